$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: split the paragraph "我is阿卡表现及保险科技进步" into three runs
# separated by grammar-check proofErr markers, without touching the
# paragraph's own pilcrow (paragraph mark) so the paragraph count/pPr of
# this paragraph stay untouched.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5r = $p5.Range
$p5Text = $d.Range($p5r.Start, $p5r.End - 1)

$xmlSplit = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>我is阿</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>卡表现及保险科技</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>进步</w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p5Text.InsertXML($xmlSplit)

# ---------------------------------------------------------------------------
# Step 2: on the paragraph "刷卡成本哈开户行" (now the 6th paragraph), drop the
# w:hint="eastAsia" attribute from its own paragraph-mark run properties
# (w:pPr/w:rPr/w:rFonts), and append two brand-new paragraphs after it:
#   - "2022年7月1日星期五" (typed as six separate runs)
#   - "心情烦躁的一天"
# All three paragraphs are produced together via a single InsertXML call
# whose target range spans from the start of paragraph 6 through the end of
# the document, which lets the supplied <w:pPr> actually replace the
# paragraph-mark formatting (a range confined to a single untouched
# paragraph leaves the existing pPr/rPr as-is).
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$p6r = $p6.Range
$tailRange = $d.Range($p6r.Start, $d.Content.End)

$xmlTail = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>刷卡成本哈开户行</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>2022年</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>7</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>月</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>1</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>日星期</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>五</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>心情烦躁的一天</w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tailRange.InsertXML($xmlTail)
